$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for existing N=100 block (row 1) ---
$ws.Range("A1").Value = "N=100"
$ws.Range("A1").Interior.Color = 65535

# --- New "N=365" section header (row 18) ---
$ws.Range("A18").Value = "N=365"
$ws.Range("A18").Interior.Color = 65535

# --- Row 19: labels ---
$ws.Range("A19").Value = "Time"
$ws.Range("G19").Value = "Time * CPUs"

# --- Row 20: CPU-count headers ---
$ws.Range("A20").Value = "M \ CPUs"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 4
$ws.Range("G20").Value = "M \ CPUs"
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 4

# --- Rows 21-23: timing data (plus a temporary row 24 used only to build
#     the shared-formula ranges the same way the original N=100 block was
#     built, then removed) ---
$ws.Range("A21").Value = 1000000
$ws.Range("A22").Value = 5000000
$ws.Range("A23").Value = 10000000
$ws.Range("A24").Value = 20000000

$ws.Range("B21").Value = 49.824
$ws.Range("B22").Value = 247.946
$ws.Range("B23").Value = 496.328
$ws.Range("B24").Value = 992.656

$ws.Range("C21").Value = 25.755
$ws.Range("C22").Value = 127.39
$ws.Range("C23").Value = 253.424
$ws.Range("C24").Value = 506.848

$ws.Range("D21").Value = 16.779
$ws.Range("D22").Value = 83.708
$ws.Range("D23").Value = 167.04
$ws.Range("D24").Value = 334.08

$ws.Range("E21").Value = 12.632
$ws.Range("E22").Value = 63.932
$ws.Range("E23").Value = 130.213
$ws.Range("E24").Value = 260.426

$ws.Range("G21").Value = 1000000
$ws.Range("G22").Value = 5000000
$ws.Range("G23").Value = 10000000
$ws.Range("G24").Value = 20000000

$ws.Range("H21").Formula = "=B21*B`$3"
$ws.Range("I21:I24").Formula = "=C21*C`$3"
$ws.Range("J21:J24").Formula = "=D21*D`$3"
$ws.Range("K21:K24").Formula = "=E21*E`$3"
$ws.Range("H22:H24").Formula = "=B22*B`$3"

$ws.Rows("24:24").Delete() | Out-Null

# --- Row 27: "Error" heading ---
$ws.Range("A27").Value = "Error"

# --- Row 28: CPU-count headers ---
$ws.Range("A28").Value = "M \ CPUs"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 4

# --- Rows 29-31: error data (written as plain decimals; this PowerShell
#     dialect does not accept scientific-notation numeric literals) ---
$ws.Range("A29").Value = 1000000
$ws.Range("B29").Value = 0.010855991
$ws.Range("C29").Value = 0.0190141
$ws.Range("D29").Value = 0.015882194
$ws.Range("E29").Value = 0.006467221

$ws.Range("A30").Value = 5000000
$ws.Range("B30").Value = 0.003155988
$ws.Range("C30").Value = 0.005549817
$ws.Range("D30").Value = 0.00514671
$ws.Range("E30").Value = 0.001827144

$ws.Range("A31").Value = 10000000
$ws.Range("B31").Value = 0.000430817
$ws.Range("C31").Value = 0.001016266
$ws.Range("D31").Value = 0.000974426
$ws.Range("E31").Value = 0.000241624

# --- Final selection state ---
$ws.Range("E25").Select() | Out-Null
